$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded between the existing rows 72
# (2021-10-04) and what was row 73 (2021-09-06, now shifted to row 74).
# Insert a fresh row at 73, pushing the old rows 73:91 down to 74:92, then
# populate the new row with the latest reading.
$ws.Rows.Item(73).Insert()

$ws.Range("A73").Value = 9
$ws.Range("B73").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C73").Value = "Metropolitana"
$ws.Range("D73").Value = 44641
$ws.Range("E73").Value = 13
$ws.Range("F73").Value = 100114007
$ws.Range("G73").Value = "Jengibre"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 610
$ws.Range("K73").Value = 14000
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = 14500
$ws.Range("N73").Value = "`$/caja 13 kilos"
$ws.Range("O73").Value = "Per" + [char]0x00FA
$ws.Range("P73").Value = 1115
$ws.Range("Q73").Value = 13
$ws.Range("R73").Value = "Hortaliza"
